$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.45425133333334
$ws.Range("H2").Value = 55.362754
$ws.Range("I2").Value = 0.3394975357727733
$ws.Range("J2").Value = 0.3394975357727733
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 11.36407333333333
$ws.Range("N2").Value = 34.09222
$ws.Range("O2").Value = 0.03922668827193482
$ws.Range("P2").Value = 0.03922668827193482
$ws.Range("Q2").Value = 209.7154654637644
$ws.Range("R2").Value = 1887.43918917388
$ws.Range("S2").Value = 0.01331736400484862
$ws.Range("T2").Value = 0.01331736400484862

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.45425133333334
$ws.Range("H3").Value = 55.362754
$ws.Range("I3").Value = 0.3394975357727733
$ws.Range("J3").Value = 0.3394975357727733
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 36.67092
$ws.Range("N3").Value = 110.01276
$ws.Range("O3").Value = 0.126581262307212
$ws.Range("P3").Value = 0.126581262307212
$ws.Range("Q3").Value = 676.73437430456
$ws.Range("R3").Value = 6090.609368741039
$ws.Range("S3").Value = 0.04297402662830551
$ws.Range("T3").Value = 0.04297402662830551

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.45425133333334
$ws.Range("H4").Value = 55.362754
$ws.Range("I4").Value = 0.3394975357727733
$ws.Range("J4").Value = 0.3394975357727733
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.688545666666666
$ws.Range("N4").Value = 14.065637
$ws.Range("O4").Value = 0.01618399617112621
$ws.Range("P4").Value = 0.01618399617112621
$ws.Range("Q4").Value = 86.52360012047755
$ws.Range("R4").Value = 778.712401084298
$ws.Range("S4").Value = 0.005494426819053347
$ws.Range("T4").Value = 0.005494426819053347

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 18.45425133333334
$ws.Range("H5").Value = 55.362754
$ws.Range("I5").Value = 0.3394975357727733
$ws.Range("J5").Value = 0.3394975357727733
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 236.9790546666667
$ws.Range("N5").Value = 710.9371639999999
$ws.Range("O5").Value = 0.8180080532497269
$ws.Range("P5").Value = 0.8180080532497269
$ws.Range("Q5").Value = 4373.271035554407
$ws.Range("R5").Value = 39359.43931998965
$ws.Range("S5").Value = 0.2777117183205658
$ws.Range("T5").Value = 0.2777117183205658

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 15.88630666666666
$ws.Range("H6").Value = 47.65891999999999
$ws.Range("I6").Value = 0.2922557988641919
$ws.Range("J6").Value = 0.292255798864192
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 11.36407333333333
$ws.Range("N6").Value = 34.09222
$ws.Range("O6").Value = 0.03922668827193482
$ws.Range("P6").Value = 0.03922668827193482
$ws.Range("Q6").Value = 180.5331539558222
$ws.Range("R6").Value = 1624.7983856024
$ws.Range("S6").Value = 0.01146422711771094
$ws.Range("T6").Value = 0.01146422711771094

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 15.88630666666666
$ws.Range("H7").Value = 47.65891999999999
$ws.Range("I7").Value = 0.2922557988641919
$ws.Range("J7").Value = 0.292255798864192
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 36.67092
$ws.Range("N7").Value = 110.01276
$ws.Range("O7").Value = 0.126581262307212
$ws.Range("P7").Value = 0.126581262307212
$ws.Range("Q7").Value = 582.5654808687998
$ws.Range("R7").Value = 5243.089327819199
$ws.Range("S7").Value = 0.03699410793683208
$ws.Range("T7").Value = 0.03699410793683208

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 15.88630666666666
$ws.Range("H8").Value = 47.65891999999999
$ws.Range("I8").Value = 0.2922557988641919
$ws.Range("J8").Value = 0.292255798864192
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.688545666666666
$ws.Range("N8").Value = 14.065637
$ws.Range("O8").Value = 0.01618399617112621
$ws.Range("P8").Value = 0.01618399617112621
$ws.Range("Q8").Value = 74.48367428133776
$ws.Range("R8").Value = 670.3530685320399
$ws.Range("S8").Value = 0.004729866729807514
$ws.Range("T8").Value = 0.004729866729807515

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 15.88630666666666
$ws.Range("H9").Value = 47.65891999999999
$ws.Range("I9").Value = 0.2922557988641919
$ws.Range("J9").Value = 0.292255798864192
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 236.9790546666667
$ws.Range("N9").Value = 710.9371639999999
$ws.Range("O9").Value = 0.8180080532497269
$ws.Range("P9").Value = 0.8180080532497269
$ws.Range("Q9").Value = 3764.72193601143
$ws.Range("R9").Value = 33882.49742410288
$ws.Range("S9").Value = 0.2390675970798414
$ws.Range("T9").Value = 0.2390675970798414

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.826026
$ws.Range("H10").Value = 14.478078
$ws.Range("I10").Value = 0.08878300750222795
$ws.Range("J10").Value = 0.08878300750222796
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 11.36407333333333
$ws.Range("N10").Value = 34.09222
$ws.Range("O10").Value = 0.03922668827193482
$ws.Range("P10").Value = 0.03922668827193482
$ws.Range("Q10").Value = 54.84331337257332
$ws.Range("R10").Value = 493.58982035316
$ws.Range("S10").Value = 0.003482663359134747
$ws.Range("T10").Value = 0.003482663359134747

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.826026
$ws.Range("H11").Value = 14.478078
$ws.Range("I11").Value = 0.08878300750222795
$ws.Range("J11").Value = 0.08878300750222796
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 36.67092
$ws.Range("N11").Value = 110.01276
$ws.Range("O11").Value = 0.126581262307212
$ws.Range("P11").Value = 0.126581262307212
$ws.Range("Q11").Value = 176.97481336392
$ws.Range("R11").Value = 1592.77332027528
$ws.Range("S11").Value = 0.01123826516106269
$ws.Range("T11").Value = 0.01123826516106269

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.826026
$ws.Range("H12").Value = 14.478078
$ws.Range("I12").Value = 0.08878300750222795
$ws.Range("J12").Value = 0.08878300750222796
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 4.688545666666666
$ws.Range("N12").Value = 14.065637
$ws.Range("O12").Value = 0.01618399617112621
$ws.Range("P12").Value = 0.01618399617112621
$ws.Range("Q12").Value = 22.62704328952066
$ws.Range("R12").Value = 203.643389605686
$ws.Range("S12").Value = 0.001436863853477127
$ws.Range("T12").Value = 0.001436863853477127

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.826026
$ws.Range("H13").Value = 14.478078
$ws.Range("I13").Value = 0.08878300750222795
$ws.Range("J13").Value = 0.08878300750222796
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 236.9790546666667
$ws.Range("N13").Value = 710.9371639999999
$ws.Range("O13").Value = 0.8180080532497269
$ws.Range("P13").Value = 0.8180080532497269
$ws.Range("Q13").Value = 1143.667079276755
$ws.Range("R13").Value = 10293.00371349079
$ws.Range("S13").Value = 0.07262521512855338
$ws.Range("T13").Value = 0.0726252151285534

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 15.19095733333334
$ws.Range("H14").Value = 45.572872
$ws.Range("I14").Value = 0.2794636578608068
$ws.Range("J14").Value = 0.2794636578608069
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 11.36407333333333
$ws.Range("N14").Value = 34.09222
$ws.Range("O14").Value = 0.03922668827193482
$ws.Range("P14").Value = 0.03922668827193482
$ws.Range("Q14").Value = 172.6311531395378
$ws.Range("R14").Value = 1553.68037825584
$ws.Range("S14").Value = 0.01096243379024052
$ws.Range("T14").Value = 0.01096243379024052

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 15.19095733333334
$ws.Range("H15").Value = 45.572872
$ws.Range("I15").Value = 0.2794636578608068
$ws.Range("J15").Value = 0.2794636578608069
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 36.67092
$ws.Range("N15").Value = 110.01276
$ws.Range("O15").Value = 0.126581262307212
$ws.Range("P15").Value = 0.126581262307212
$ws.Range("Q15").Value = 557.06638109408
$ws.Range("R15").Value = 5013.59742984672
$ws.Range("S15").Value = 0.03537486258101175
$ws.Range("T15").Value = 0.03537486258101176

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 15.19095733333334
$ws.Range("H16").Value = 45.572872
$ws.Range("I16").Value = 0.2794636578608068
$ws.Range("J16").Value = 0.2794636578608069
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 4.688545666666666
$ws.Range("N16").Value = 14.065637
$ws.Range("O16").Value = 0.01618399617112621
$ws.Range("P16").Value = 0.01618399617112621
$ws.Range("Q16").Value = 71.22349717771823
$ws.Range("R16").Value = 641.011474599464
$ws.Range("S16").Value = 0.004522838768788224
$ws.Range("T16").Value = 0.004522838768788225

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 15.19095733333334
$ws.Range("H17").Value = 45.572872
$ws.Range("I17").Value = 0.2794636578608068
$ws.Range("J17").Value = 0.2794636578608069
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 236.9790546666667
$ws.Range("N17").Value = 710.9371639999999
$ws.Range("O17").Value = 0.8180080532497269
$ws.Range("P17").Value = 0.8180080532497269
$ws.Range("Q17").Value = 3599.938708335001
$ws.Range("R17").Value = 32399.44837501501
$ws.Range("S17").Value = 0.2286035227207663
$ws.Range("T17").Value = 0.2286035227207664
